# Productivity adjustment set to zero
# Sheet "F2C Jobs" column J ("Productivity") values are zeroed out.
# A handful of rows (9,10,11,12,13,26,27,28) had their Productivity
# (and, for 11-13, the whole I:M block) cells cleared entirely rather
# than set to 0, matching the upstream edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("F2C Jobs")

# Rows whose J (Productivity) value becomes 0
$zeroRows = @(2,3,4,5,6,7,8) + @(14..25) + @(32..58)
foreach ($r in $zeroRows) {
    $ws.Cells.Item($r, 10).Value = 0
}

# Rows 9 and 10: Productivity cell cleared entirely (no value at all)
foreach ($r in @(9,10)) {
    $ws.Cells.Item($r, 10).ClearContents()
}

# Rows 26,27,28: both Uncertainty Range (I) and Productivity (J) cells cleared
foreach ($r in @(26,27,28)) {
    $ws.Range("I$r`:J$r").ClearContents()
}

# Rows 11,12,13: entire I:M block cleared (placeholder cells removed)
foreach ($r in @(11,12,13)) {
    $ws.Range("I$r`:M$r").ClearContents()
}
